# Add files via upload
# - Bold/size-12 header row (with a taller row height)
# - A new data row (row 3) with a date and two descriptive notes
# - Portrait page setup

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (A1:D1): bold, 12pt, taller row ---
$ws.Range("A1:D1").Font.Bold = $true
$ws.Range("A1:D1").Font.Size = 12
$ws.Rows.Item(1).RowHeight = 15.75

# --- New row 3: reuse row 2's date style for A3, then fill in the text ---
$ws.Range("A2").Copy()
$ws.Range("A3").PasteSpecial(-4122)
$ws.Range("A3").Value = [DateTime]"2019-12-06"

$ws.Range("C3").Value = "feasibility study on this project"
$ws.Range("B3").Value = "Concepts of how the login page in that using the forgot and OTP sending the email account"

# --- Page setup ---
$ws.PageSetup.Orientation = 1

# --- Selection matches the author's last active cell ---
$ws.Range("C3").Select() | Out-Null
